# Refresh the Pspn-Ret ligand-receptor table with updated TPM-derived NATMI output.
# The sending/target cluster set grew (Inflammatory-Mac, MuSCs, Resolving-Mac now all
# appear as Sending cluster; ECs/FAPs/MuSCs as Target cluster), so the table grows
# from 8 data rows to 12 (rows 2-13) with refreshed numeric columns throughout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data block (row 1 header is untouched) before writing the refreshed table
$ws.Range("A2:T9").Clear()

# Row 2
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Pspn"
$ws.Cells.Item(2,3).Value = "Ret"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 3.042572666666667
$ws.Cells.Item(2,8).Value = 9.127718
$ws.Cells.Item(2,9).Value = 0.7055529458943673
$ws.Cells.Item(2,10).Value = 0.7055529458943673
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 3.21276
$ws.Cells.Item(2,14).Value = 9.63828
$ws.Cells.Item(2,15).Value = 0.3952107490920524
$ws.Cells.Item(2,16).Value = 0.3952107490920524
$ws.Cells.Item(2,17).Value = 9.775055760559999
$ws.Cells.Item(2,18).Value = 87.97550184504
$ws.Cells.Item(2,19).Value = 0.2788421082710173
$ws.Cells.Item(2,20).Value = 0.2788421082710173

# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Pspn"
$ws.Cells.Item(3,3).Value = "Ret"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 3.042572666666667
$ws.Cells.Item(3,8).Value = 9.127718
$ws.Cells.Item(3,9).Value = 0.7055529458943673
$ws.Cells.Item(3,10).Value = 0.7055529458943673
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 3.884996
$ws.Cells.Item(3,14).Value = 11.654988
$ws.Cells.Item(3,15).Value = 0.4779044122124365
$ws.Cells.Item(3,16).Value = 0.4779044122124365
$ws.Cells.Item(3,17).Value = 11.82038263970933
$ws.Cells.Item(3,18).Value = 106.383443757384
$ws.Cells.Item(3,19).Value = 0.3371868658924006
$ws.Cells.Item(3,20).Value = 0.3371868658924006

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Pspn"
$ws.Cells.Item(4,3).Value = "Ret"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 3.042572666666667
$ws.Cells.Item(4,8).Value = 9.127718
$ws.Cells.Item(4,9).Value = 0.7055529458943673
$ws.Cells.Item(4,10).Value = 0.7055529458943673
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.031476333333333
$ws.Cells.Item(4,14).Value = 3.094429
$ws.Cells.Item(4,15).Value = 0.1268848386955111
$ws.Cells.Item(4,16).Value = 0.1268848386955111
$ws.Cells.Item(4,17).Value = 3.138341698113556
$ws.Cells.Item(4,18).Value = 28.245075283022
$ws.Cells.Item(4,19).Value = 0.08952397173094949
$ws.Cells.Item(4,20).Value = 0.08952397173094949

# Row 5
$ws.Cells.Item(5,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(5,2).Value = "Pspn"
$ws.Cells.Item(5,3).Value = "Ret"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.271593
$ws.Cells.Item(5,8).Value = 0.814779
$ws.Cells.Item(5,9).Value = 0.06298066216581918
$ws.Cells.Item(5,10).Value = 0.06298066216581918
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 3.21276
$ws.Cells.Item(5,14).Value = 9.63828
$ws.Cells.Item(5,15).Value = 0.3952107490920524
$ws.Cells.Item(5,16).Value = 0.3952107490920524
$ws.Cells.Item(5,17).Value = 0.87256312668
$ws.Cells.Item(5,18).Value = 7.85306814012
$ws.Cells.Item(5,19).Value = 0.02489063467286689
$ws.Cells.Item(5,20).Value = 0.02489063467286689

# Row 6
$ws.Cells.Item(6,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(6,2).Value = "Pspn"
$ws.Cells.Item(6,3).Value = "Ret"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.271593
$ws.Cells.Item(6,8).Value = 0.814779
$ws.Cells.Item(6,9).Value = 0.06298066216581918
$ws.Cells.Item(6,10).Value = 0.06298066216581918
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 3.884996
$ws.Cells.Item(6,14).Value = 11.654988
$ws.Cells.Item(6,15).Value = 0.4779044122124365
$ws.Cells.Item(6,16).Value = 0.4779044122124365
$ws.Cells.Item(6,17).Value = 1.055137718628
$ws.Cells.Item(6,18).Value = 9.496239467652
$ws.Cells.Item(6,19).Value = 0.03009873633310585
$ws.Cells.Item(6,20).Value = 0.03009873633310585

# Row 7
$ws.Cells.Item(7,1).Value = "Inflammatory-Mac"
$ws.Cells.Item(7,2).Value = "Pspn"
$ws.Cells.Item(7,3).Value = "Ret"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.271593
$ws.Cells.Item(7,8).Value = 0.814779
$ws.Cells.Item(7,9).Value = 0.06298066216581918
$ws.Cells.Item(7,10).Value = 0.06298066216581918
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.031476333333333
$ws.Cells.Item(7,14).Value = 3.094429
$ws.Cells.Item(7,15).Value = 0.1268848386955111
$ws.Cells.Item(7,16).Value = 0.1268848386955111
$ws.Cells.Item(7,17).Value = 0.2801417517990001
$ws.Cells.Item(7,18).Value = 2.521275766191
$ws.Cells.Item(7,19).Value = 0.007991291159846448
$ws.Cells.Item(7,20).Value = 0.007991291159846448

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Pspn"
$ws.Cells.Item(8,3).Value = "Ret"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = 0.3333333333333333
$ws.Cells.Item(8,7).Value = 0.3383276666666666
$ws.Cells.Item(8,8).Value = 1.014983
$ws.Cells.Item(8,9).Value = 0.07845600024920825
$ws.Cells.Item(8,10).Value = 0.07845600024920826
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 3.21276
$ws.Cells.Item(8,14).Value = 9.63828
$ws.Cells.Item(8,15).Value = 0.3952107490920524
$ws.Cells.Item(8,16).Value = 0.3952107490920524
$ws.Cells.Item(8,17).Value = 1.08696559436
$ws.Cells.Item(8,18).Value = 9.782690349239999
$ws.Cells.Item(8,19).Value = 0.03100665462925585
$ws.Cells.Item(8,20).Value = 0.03100665462925585

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Pspn"
$ws.Cells.Item(9,3).Value = "Ret"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 0.3333333333333333
$ws.Cells.Item(9,7).Value = 0.3383276666666666
$ws.Cells.Item(9,8).Value = 1.014983
$ws.Cells.Item(9,9).Value = 0.07845600024920825
$ws.Cells.Item(9,10).Value = 0.07845600024920826
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 3.884996
$ws.Cells.Item(9,14).Value = 11.654988
$ws.Cells.Item(9,15).Value = 0.4779044122124365
$ws.Cells.Item(9,16).Value = 0.4779044122124365
$ws.Cells.Item(9,17).Value = 1.314401631689333
$ws.Cells.Item(9,18).Value = 11.829614685204
$ws.Cells.Item(9,19).Value = 0.03749446868363664
$ws.Cells.Item(9,20).Value = 0.03749446868363664

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Pspn"
$ws.Cells.Item(10,3).Value = "Ret"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0.3333333333333333
$ws.Cells.Item(10,7).Value = 0.3383276666666666
$ws.Cells.Item(10,8).Value = 1.014983
$ws.Cells.Item(10,9).Value = 0.07845600024920825
$ws.Cells.Item(10,10).Value = 0.07845600024920826
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.031476333333333
$ws.Cells.Item(10,14).Value = 3.094429
$ws.Cells.Item(10,15).Value = 0.1268848386955111
$ws.Cells.Item(10,16).Value = 0.1268848386955111
$ws.Cells.Item(10,17).Value = 0.3489769810785556
$ws.Cells.Item(10,18).Value = 3.140792829707
$ws.Cells.Item(10,19).Value = 0.009954876936315769
$ws.Cells.Item(10,20).Value = 0.009954876936315771

# Row 11
$ws.Cells.Item(11,1).Value = "Resolving-Mac"
$ws.Cells.Item(11,2).Value = "Pspn"
$ws.Cells.Item(11,3).Value = "Ret"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 0.3333333333333333
$ws.Cells.Item(11,7).Value = 0.6598303333333333
$ws.Cells.Item(11,8).Value = 1.979491
$ws.Cells.Item(11,9).Value = 0.1530103916906052
$ws.Cells.Item(11,10).Value = 0.1530103916906052
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 3.21276
$ws.Cells.Item(11,14).Value = 9.63828
$ws.Cells.Item(11,15).Value = 0.3952107490920524
$ws.Cells.Item(11,16).Value = 0.3952107490920524
$ws.Cells.Item(11,17).Value = 2.11987650172
$ws.Cells.Item(11,18).Value = 19.07888851548
$ws.Cells.Item(11,19).Value = 0.06047135151891242
$ws.Cells.Item(11,20).Value = 0.06047135151891242

# Row 12
$ws.Cells.Item(12,1).Value = "Resolving-Mac"
$ws.Cells.Item(12,2).Value = "Pspn"
$ws.Cells.Item(12,3).Value = "Ret"
$ws.Cells.Item(12,4).Value = "FAPs"
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(12,6).Value = 0.3333333333333333
$ws.Cells.Item(12,7).Value = 0.6598303333333333
$ws.Cells.Item(12,8).Value = 1.979491
$ws.Cells.Item(12,9).Value = 0.1530103916906052
$ws.Cells.Item(12,10).Value = 0.1530103916906052
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 3.884996
$ws.Cells.Item(12,14).Value = 11.654988
$ws.Cells.Item(12,15).Value = 0.4779044122124365
$ws.Cells.Item(12,16).Value = 0.4779044122124365
$ws.Cells.Item(12,17).Value = 2.563438205678666
$ws.Cells.Item(12,18).Value = 23.070943851108
$ws.Cells.Item(12,19).Value = 0.07312434130329333
$ws.Cells.Item(12,20).Value = 0.07312434130329333

# Row 13
$ws.Cells.Item(13,1).Value = "Resolving-Mac"
$ws.Cells.Item(13,2).Value = "Pspn"
$ws.Cells.Item(13,3).Value = "Ret"
$ws.Cells.Item(13,4).Value = "MuSCs"
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 0.3333333333333333
$ws.Cells.Item(13,7).Value = 0.6598303333333333
$ws.Cells.Item(13,8).Value = 1.979491
$ws.Cells.Item(13,9).Value = 0.1530103916906052
$ws.Cells.Item(13,10).Value = 0.1530103916906052
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 1.031476333333333
$ws.Cells.Item(13,14).Value = 3.094429
$ws.Cells.Item(13,15).Value = 0.1268848386955111
$ws.Cells.Item(13,16).Value = 0.1268848386955111
$ws.Cells.Item(13,17).Value = 0.6805993728487778
$ws.Cells.Item(13,18).Value = 6.125394355639
$ws.Cells.Item(13,19).Value = 0.01941469886839941
$ws.Cells.Item(13,20).Value = 0.01941469886839941
